$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H (row 1), matching header style of existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean values for rows 2-14 in columns F, G, H -- default False, except H6 = True
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

$ws.Cells.Item(6, 8).Value = $true
